$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Change C6 value from YES to NO
$ws.Range("C6").Value = "NO"

# Highlight row 5 (A5:C5) with yellow fill by selecting the entire row
$ws.Rows("5:5").Interior.Color = 65535

# Update selection to C7
$ws.Range("C7").Select()
